$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# ---------------------------------------------------------------------------
# 1) Row 20 <-> Row 22 : swap match details (same matchday, E column unchanged)
# ---------------------------------------------------------------------------
$v20 = Get-RowValues 20
$v22 = Get-RowValues 22
Set-RowValues 20 $v22
Set-RowValues 22 $v20

# ---------------------------------------------------------------------------
# 2) Rows 27, 28, 29 : 3-way rotation (new27=old28, new28=old29, new29=old27)
# ---------------------------------------------------------------------------
$v27 = Get-RowValues 27
$v28 = Get-RowValues 28
$v29 = Get-RowValues 29
Set-RowValues 27 $v28
Set-RowValues 28 $v29
Set-RowValues 29 $v27

# ---------------------------------------------------------------------------
# 3) Rows 41, 43, 44 : 3-way rotation (new41=old43, new43=old44, new44=old41)
# ---------------------------------------------------------------------------
$v41 = Get-RowValues 41
$v43 = Get-RowValues 43
$v44 = Get-RowValues 44
Set-RowValues 41 $v43
Set-RowValues 43 $v44
Set-RowValues 44 $v41

# ---------------------------------------------------------------------------
# 4) Row 47 <-> Row 48 : swap match details
# ---------------------------------------------------------------------------
$v47 = Get-RowValues 47
$v48 = Get-RowValues 48
Set-RowValues 47 $v48
Set-RowValues 48 $v47

# ---------------------------------------------------------------------------
# 5) Append 4 new match rows (59-62), cloning the formatting of row 58
# ---------------------------------------------------------------------------
$ws.Range("A58:V58").Copy()
$ws.Range("A59:V62").PasteSpecial(-4122)

$newRows = @{
    59 = @{ A=58; B="poland"; C="iii-liga-group-iii"; D="2023-2024"; E=45192.45833333334;
            F="Bytom Odrzanski"; G=2; H="Rakow II"; I=1; J=2.4; K="22/09/2023 22:42";
            L=1.85; M="23/09/2023 10:46"; N=3.4; O="22/09/2023 22:42"; P=3.78; Q="23/09/2023 10:46";
            R=2.5; S="22/09/2023 22:42"; T=3.3; U="23/09/2023 10:46";
            V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/bytom-odrzanski-rks-rakow-czestochowa/Cdunz9B8/" };
    60 = @{ A=59; B="poland"; C="iii-liga-group-iii"; D="2023-2024"; E=45192.5;
            F="Carina Gubin"; G=1; H="Kluczbork"; I=1; J=2.44; K="21/09/2023 23:12";
            L=2.52; M="23/09/2023 11:58"; N=3.27; O="21/09/2023 23:12"; P=3.43; Q="23/09/2023 11:58";
            R=2.32; S="21/09/2023 23:12"; T=2.4; U="23/09/2023 11:51";
            V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/carina-gubin-kluczbork/UDamF7Zr/" };
    61 = @{ A=60; B="poland"; C="iii-liga-group-iii"; D="2023-2024"; E=45192.54166666666;
            F="Sleza Wroclaw"; G=2; H="Bielsko-Biala"; I=2; J=2.72; K="22/09/2023 00:13";
            L=2.5; M="23/09/2023 12:49"; N=3.37; O="22/09/2023 00:13"; P=3.57; Q="23/09/2023 12:44";
            R=2.07; S="22/09/2023 00:13"; T=2.36; U="23/09/2023 12:49";
            V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/sleza-wroclaw-rekord-bielsko-biala/OzaiERlk/" };
    62 = @{ A=61; B="poland"; C="iii-liga-group-iii"; D="2023-2024"; E=45192.54166666666;
            F="Starowice Dolne"; G=0; H="Stilon Gorzow"; I=7; J=1.94; K="22/09/2023 00:13";
            L=2.17; M="23/09/2023 11:44"; N=3.47; O="22/09/2023 00:13"; P=3.53; Q="23/09/2023 11:44";
            R=2.99; S="22/09/2023 00:13"; T=2.77; U="23/09/2023 11:44";
            V="https://www.betexplorer.com/football/poland/iii-liga-group-iii/starowice-dolne-stilon-gorzow/K2h0C5J1/" };
}

$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($r in @(59,60,61,62)) {
    $rowData = $newRows[$r]
    foreach ($col in $allCols) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}

"Edit complete"
